$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.700.28"
Set-TextValue $ws.Range("E2") "  -1.99%  "

Set-TextValue $ws.Range("D3") "3.390.53"
Set-TextValue $ws.Range("E3") "  -2.34%  "

Set-TextValue $ws.Range("E4") "  -0.14%  "

Set-TextValue $ws.Range("D5") "405.30"
Set-TextValue $ws.Range("E5") "  -2.52%  "

Set-TextValue $ws.Range("D6") "133.64"
Set-TextValue $ws.Range("E6") "  +8.64%  "

Set-TextValue $ws.Range("E7") "  -1.51%  "

Set-TextValue $ws.Range("E8") "  -0.06%  "

Set-TextValue $ws.Range("D9") "0.670"
Set-TextValue $ws.Range("E9") "  -3.13%  "

Set-TextValue $ws.Range("E10") "  -8.10%  "

Set-TextValue $ws.Range("D11") "42.47"
Set-TextValue $ws.Range("E11") "  +2.38%  "

Set-TextValue $ws.Range("E12") "  -1.35%  "

Set-TextValue $ws.Range("D13") "3.913.01"
Set-TextValue $ws.Range("E13") "  -2.99%  "

Set-TextValue $ws.Range("D14") "8.43"
Set-TextValue $ws.Range("E14") "  -2.42%  "

Set-TextValue $ws.Range("E15") "  -1.31%  "

Set-TextValue $ws.Range("D16") "3.392.44"
Set-TextValue $ws.Range("E16") "  -2.75%  "

Set-TextValue $ws.Range("D17") "61.591.93"
Set-TextValue $ws.Range("E17") "  -2.14%  "

Set-TextValue $ws.Range("E18") "  -1.76%  "

Set-TextValue $ws.Range("D19") "11.02"
Set-TextValue $ws.Range("E19") "  +0.66%  "

Set-TextValue $ws.Range("E20") "  -8.37%  "

Set-TextValue $ws.Range("E21") "  -3.48%  "

Set-TextValue $ws.Range("D22") "85.49"
Set-TextValue $ws.Range("E22") "  +3.93%  "

Set-TextValue $ws.Range("D23") "314.92"
Set-TextValue $ws.Range("E23") "  -0.90%  "

Set-TextValue $ws.Range("D24") "12.74"
Set-TextValue $ws.Range("E24") "  -1.21%  "

Set-TextValue $ws.Range("E25") "  -1.85%  "

Set-TextValue $ws.Range("E26") "  +11.18%  "

Set-TextValue $ws.Range("D27") "29.56"
Set-TextValue $ws.Range("E27") "  -4.86%  "

Set-TextValue $ws.Range("D28") "8.31"
Set-TextValue $ws.Range("E28") "  +5.48%  "

Set-TextValue $ws.Range("D29") "7.68"
Set-TextValue $ws.Range("E29") "  -1.65%  "

Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.66"
Set-TextValue $ws.Range("E30") "  +4.05%  "

Set-TextValue $ws.Range("E31") "  -1.97%  "

Set-TextValue $ws.Range("B32") "Hedera"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D32") "0.116"
Set-TextValue $ws.Range("E32") "  -0.90%  "

Set-TextValue $ws.Range("B33") "Dai"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D33") "1.00"
Set-TextValue $ws.Range("E33") "  -0.62%  "

Set-TextValue $ws.Range("B34") "Cosmos"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D34") "11.33"
Set-TextValue $ws.Range("E34") "  -2.40%  "

Set-TextValue $ws.Range("D35") "41.45"
Set-TextValue $ws.Range("E35") "  -1.53%  "

Set-TextValue $ws.Range("D36") "0.0480"
Set-TextValue $ws.Range("E36") "  -2.02%  "

Set-TextValue $ws.Range("D37") "51.81"
Set-TextValue $ws.Range("E37") "  -0.39%  "

Set-TextValue $ws.Range("E38") "  -0.07%  "

Set-TextValue $ws.Range("E39") "  -1.88%  "

Set-TextValue $ws.Range("D40") "2.94"
Set-TextValue $ws.Range("E40") "  -2.87%  "

Set-TextValue $ws.Range("D41") "139.25"
Set-TextValue $ws.Range("E41") "  +2.47%  "

Set-TextValue $ws.Range("E42") "  -0.83%  "

Set-TextValue $ws.Range("E43") "  -1.20%  "

Set-TextValue $ws.Range("D44") "0.296"
Set-TextValue $ws.Range("E44") "  +4.74%  "

Set-TextValue $ws.Range("D45") "3.97"
Set-TextValue $ws.Range("E45") "  +2.21%  "

Set-TextValue $ws.Range("D46") "16.62"
Set-TextValue $ws.Range("E46") "  -1.42%  "

Set-TextValue $ws.Range("E47") "  -1.73%  "

Set-TextValue $ws.Range("D48") "21.41"
Set-TextValue $ws.Range("E48") "  -2.34%  "

Set-TextValue $ws.Range("D49") "2.121.57"
Set-TextValue $ws.Range("E49") "  -2.98%  "

Set-TextValue $ws.Range("D50") "2.30"
Set-TextValue $ws.Range("E50") "  -6.45%  "

Set-TextValue $ws.Range("E51") "  +0.87%  "
